# Auto-committed on 2023/09/22 週五 16:59:26.34
# Adds two new rows (CustIdErrFg / SpouseIdErrFg) to the DBD layout sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

# --- Row 84 : 身份證字號/統一編號錯誤註記 -----------------------------------
$ws.Range("A84").Formula = "=A83+1"
$ws.Range("B84").Value = "CustIdErrFg"
$ws.Range("C84").Value = "身份證字號/統一編號錯誤註記"
$ws.Range("D84").Value = "VARCHAR2"
$ws.Range("E84").Value = 1
$ws.Range("G84").Value = "Y:`nA:舊資料轉換`nB:舊資料轉換"
$ws.Range("H84").Value = 45190

# --- Row 85 : 配偶身份證號/負責人身分證錯誤註記 -----------------------------
$ws.Range("A85").Formula = "=A84+1"
$ws.Range("B85").Value = "SpouseIdErrFg"
$ws.Range("C85").Value = "配偶身份證號/負責人身分證錯誤註記"
$ws.Range("D85").Value = "VARCHAR2"
$ws.Range("E85").Value = 1
$ws.Range("G85").Value = "Y:`nA:舊資料轉換`nB:舊資料轉換"
$ws.Range("H85").Value = 45190

# --- Formatting: reuse the workbook's existing "new row" styles ------------
# Columns D/E use the style already applied at E16 (red 標楷體, wrap, centered)
$ws.Range("E16").Copy()
[void]$ws.Range("D84:E84").PasteSpecial(-4122)
[void]$ws.Range("D85:E85").PasteSpecial(-4122)

# Columns A/B/C/F/G use the style already applied at F16 (red 標楷體, wrap, top)
$ws.Range("F16").Copy()
[void]$ws.Range("A84").PasteSpecial(-4122)
[void]$ws.Range("B84").PasteSpecial(-4122)
[void]$ws.Range("C84").PasteSpecial(-4122)
[void]$ws.Range("F84").PasteSpecial(-4122)
[void]$ws.Range("G84").PasteSpecial(-4122)
[void]$ws.Range("A85").PasteSpecial(-4122)
[void]$ws.Range("B85").PasteSpecial(-4122)
[void]$ws.Range("C85").PasteSpecial(-4122)
[void]$ws.Range("F85").PasteSpecial(-4122)
[void]$ws.Range("G85").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Date column H: build the short-date style once on H84, then clone it onto
# H85 via PasteSpecial so both cells share the same cellXf (no duplicate style).
$ws.Range("H84").NumberFormat = "mm-dd-yy"
$ws.Range("H84").Copy()
[void]$ws.Range("H85").PasteSpecial(-4122)
$ws.Range("H85").Value = 45190
$excel.CutCopyMode = $false

# Row heights grow to fit the wrapped three-line note in column G
$ws.Rows.Item(84).RowHeight = 48.6
$ws.Rows.Item(85).RowHeight = 48.6

# --- View state: scroll down and select the last edited cell ---------------
[void]$ws.Range("H85").Select()
